$wb = $excel.ActiveWorkbook

# --- ALC!row 64 (diff hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 4322.5
$ws.Range("I64").Value = 3298.8
$ws.Range("J64").Value = 6028.6665
$ws.Range("K64").Value = 3298.8
$ws.Range("L64").Value = 6028.6665
$ws.Range("M64").Value = -3050.8
$ws.Range("N64").Value = -6524.6665

# --- ALC!row 67 (diff hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 4322.5
$ws.Range("I67").Value = 3298.8
$ws.Range("J67").Value = 6028.6665
$ws.Range("K67").Value = 3298.8
$ws.Range("L67").Value = 6028.6665
$ws.Range("M67").Value = -2440.8
$ws.Range("N67").Value = -7744.6665

# --- ALC!row 107 (diff hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1852716.9
$ws.Range("I107").Value = 2223100.2
$ws.Range("K107").Value = 2223100.2
$ws.Range("M107").Value = -2221180.2

# --- ALC!row 132 (diff hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 422573.47
$ws.Range("I132").Value = 486045.1
$ws.Range("J132").Value = 25876
$ws.Range("K132").Value = 1458135.3
$ws.Range("L132").Value = 77628
$ws.Range("M132").Value = -1455605.3
$ws.Range("N132").Value = -82688

# --- ARM!row 45 (diff hunk 4) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2078.25
$ws.Range("I45").Value = 1722.4
$ws.Range("J45").Value = 2671.3333
$ws.Range("K45").Value = 1722.4
$ws.Range("L45").Value = 2671.3333
$ws.Range("M45").Value = -1345.4
$ws.Range("N45").Value = -3425.3333

# --- ARM!row 61 (diff hunk 5) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2143.4314
$ws.Range("I61").Value = 1712.4222
$ws.Range("J61").Value = 5376
$ws.Range("K61").Value = 1712.4222
$ws.Range("L61").Value = 5376
$ws.Range("M61").Value = -1500.4222
$ws.Range("N61").Value = -5800

# --- ARM!row 136 (diff hunk 6) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2143.4314
$ws.Range("I136").Value = 1712.4222
$ws.Range("J136").Value = 5376
$ws.Range("K136").Value = 5137.2666
$ws.Range("L136").Value = 16128
$ws.Range("M136").Value = -2587.2666
$ws.Range("N136").Value = -21228

# --- BSM!row 15 (diff hunk 7) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

# --- BSM!row 20 (diff hunk 8) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1789.6666
$ws.Range("I20").Value = 1613.0714
$ws.Range("J20").Value = 2142.8572
$ws.Range("K20").Value = 1613.0714
$ws.Range("L20").Value = 2142.8572
$ws.Range("M20").Value = -1366.0714
$ws.Range("N20").Value = -2636.8572

# --- BSM!row 134 (diff hunk 9) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3464.7354
$ws.Range("I134").Value = 2293.0454
$ws.Range("J134").Value = 5612.8335
$ws.Range("K134").Value = 6879.1362
$ws.Range("L134").Value = 16838.5005
$ws.Range("M134").Value = -4344.1362
$ws.Range("N134").Value = -21908.5005

# --- CRP!row 31 (diff hunk 10) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3890.5652
$ws.Range("I31").Value = 1074.4117
$ws.Range("J31").Value = 11869.667
$ws.Range("K31").Value = 1074.4117
$ws.Range("L31").Value = 11869.667
$ws.Range("M31").Value = -779.4117000000001
$ws.Range("N31").Value = -12459.667

# --- CRP!row 34 (diff hunk 11) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3890.5652
$ws.Range("I34").Value = 1074.4117
$ws.Range("J34").Value = 11869.667
$ws.Range("K34").Value = 1074.4117
$ws.Range("L34").Value = 11869.667
$ws.Range("M34").Value = -872.4117000000001
$ws.Range("N34").Value = -12273.667

# --- CRP!row 105 (diff hunk 12) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1074.4286
$ws.Range("I105").Value = 942
$ws.Range("J105").Value = 1405.5
$ws.Range("K105").Value = 942
$ws.Range("L105").Value = 1405.5
$ws.Range("M105").Value = 805
$ws.Range("N105").Value = -4899.5

# --- CUL!row 4 (diff hunk 13) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30064.1
$ws.Range("I4").Value = 71.22221999999999
$ws.Range("J4").Value = 300000
$ws.Range("K4").Value = 213.66666
$ws.Range("L4").Value = 900000
$ws.Range("M4").Value = -101.66666
$ws.Range("N4").Value = -900224

# --- CUL!row 5 (diff hunk 14) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1128.6389
$ws.Range("I5").Value = 521.5454999999999
$ws.Range("J5").Value = 2082.6428
$ws.Range("K5").Value = 1564.6365
$ws.Range("L5").Value = 6247.928400000001
$ws.Range("M5").Value = -1452.6365
$ws.Range("N5").Value = -6471.928400000001

# --- CUL!row 17 (diff hunk 15) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 298.30768
$ws.Range("J17").Value = 334.36365
$ws.Range("L17").Value = 1003.09095
$ws.Range("N17").Value = -1341.09095

# --- CUL!row 98 (diff hunk 16) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 442.85715
$ws.Range("I98").Value = 440
$ws.Range("J98").Value = 450
$ws.Range("K98").Value = 1320
$ws.Range("L98").Value = 1350
$ws.Range("M98").Value = 178
$ws.Range("N98").Value = -4346

# --- CUL!row 113 (diff hunk 17) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 710.63336
$ws.Range("I113").Value = 698.0625
$ws.Range("J113").Value = 725
$ws.Range("K113").Value = 2094.1875
$ws.Range("L113").Value = 2175
$ws.Range("M113").Value = 75.8125
$ws.Range("N113").Value = -6515

# --- CUL!row 129 (diff hunk 18) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1320
$ws.Range("J129").Value = 2000
$ws.Range("L129").Value = 6000
$ws.Range("N129").Value = -16000

# --- CUL!row 132 (diff hunk 19) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1213.25
$ws.Range("I132").Value = 719.5
$ws.Range("J132").Value = 1460.125
$ws.Range("K132").Value = 6475.5
$ws.Range("L132").Value = 13141.125
$ws.Range("M132").Value = -3945.5
$ws.Range("N132").Value = -18201.125

# --- CUL!row 135 (diff hunk 20) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1128.6389
$ws.Range("I135").Value = 521.5454999999999
$ws.Range("J135").Value = 2082.6428
$ws.Range("K135").Value = 4693.9095
$ws.Range("L135").Value = 18743.7852
$ws.Range("M135").Value = -2158.9095
$ws.Range("N135").Value = -23813.7852

# --- GSM!row 107 (diff hunk 21) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1362.75
$ws.Range("I107").Value = 2249.75
$ws.Range("J107").Value = 475.75
$ws.Range("K107").Value = 2249.75
$ws.Range("L107").Value = 475.75
$ws.Range("M107").Value = -329.75
$ws.Range("N107").Value = -4315.75

# --- LTW!row 14 (diff hunk 22) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 26000
$ws.Range("I14").Value = 26000
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 26000
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -25828
$ws.Range("N14").ClearContents()

# --- LTW!row 132 (diff hunk 23) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4063.8572
$ws.Range("I132").Value = 2401
$ws.Range("J132").Value = 6281
$ws.Range("K132").Value = 7203
$ws.Range("L132").Value = 18843
$ws.Range("M132").Value = -4673
$ws.Range("N132").Value = -23903

# --- LTW!row 136 (diff hunk 24) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5257.39
$ws.Range("I136").Value = 2926.1614
$ws.Range("J136").Value = 12484.2
$ws.Range("K136").Value = 8778.484199999999
$ws.Range("L136").Value = 37452.60000000001
$ws.Range("M136").Value = -6228.484199999999
$ws.Range("N136").Value = -42552.60000000001

# --- WVR!row 19 (diff hunk 25) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 40000
$ws.Range("J19").Value = 30000
$ws.Range("L19").Value = 30000
$ws.Range("N19").Value = -30348

# --- WVR!row 92 (diff hunk 26) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
